# Update the data values in row 2 (D2, F2, H2) and move the active
# selection from D4 to C2, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 4
$ws.Range("F2").Value = -3
$ws.Range("H2").Value = 46

$ws.Range("C2").Select()
